# Update cryptocurrency price/volume data per the Nov 7 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '75.186.62'
$ws.Range("E2").Value = '  +1.29%  '

# Row 3
$ws.Range("D3").Value = '2.809.13'
$ws.Range("E3").Value = '  +6.31%  '

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").Value = "'188.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '

# Row 6
$ws.Range("D6").Value = "'593.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.00%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = "'0.545"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.90%  '

# Row 9
$ws.Range("D9").Value = "'0.191"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.06%  '

# Row 10
$ws.Range("D10").Value = '2.806.75'
$ws.Range("E10").Value = '  +6.21%  '

# Row 11
$ws.Range("D11").Value = "'0.377"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.84%  '

# Row 12
$ws.Range("E12").Value = '  -2.06%  '

# Row 13
$ws.Range("D13").Value = "'4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.07%  '

# Row 14
$ws.Range("D14").Value = '3.317.04'
$ws.Range("E14").Value = '  +6.55%  '

# Row 15
$ws.Range("D15").Value = '75.070.67'
$ws.Range("E15").Value = '  +1.27%  '

# Row 16
$ws.Range("D16").Value = "'0.0000187"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.70%  '

# Row 17
$ws.Range("D17").Value = "'26.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.65%  '

# Row 18
$ws.Range("D18").Value = '2.796.96'
$ws.Range("E18").Value = '  +5.83%  '

# Row 19
$ws.Range("D19").Value = "'9.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.85%  '

# Row 20
$ws.Range("D20").Value = "'12.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.94%  '

# Row 21
$ws.Range("D21").Value = "'377.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.46%  '

# Row 22
$ws.Range("D22").Value = "'2.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.24%  '

# Row 23
$ws.Range("D23").Value = "'4.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.02%  '

# Row 24
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("D25").Value = "'70.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.55%  '

# Row 26
$ws.Range("D26").Value = '2.950.57'
$ws.Range("E26").Value = '  +7.00%  '

# Row 27
$ws.Range("D27").Value = "'4.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.93%  '

# Row 28
$ws.Range("D28").Value = "'9.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.40%  '

# Row 29
$ws.Range("D29").Value = "'0.0000103"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.06%  '

# Row 30
$ws.Range("D30").Value = "'0.996"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.19%  '

# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = "'514.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.86%  '

# Row 32
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = "'1.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.64%  '

# Row 33
$ws.Range("D33").Value = "'7.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.16%  '

# Row 34
$ws.Range("D34").Value = "'1.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.10%  '

# Row 35
$ws.Range("E35").Value = '  +0.09%  '

# Row 36
$ws.Range("D36").Value = "'164.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.18%  '

# Row 37
$ws.Range("D37").Value = "'19.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.24%  '

# Row 38
$ws.Range("E38").Value = '  -0.05%  '

# Row 39
$ws.Range("D39").Value = "'19.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.55%  '

# Row 40
$ws.Range("E40").Value = '  -0.02%  '

# Row 41
$ws.Range("D41").Value = "'180.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.81%  '

# Row 42
$ws.Range("D42").Value = "'0.341"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.01%  '

# Row 43
$ws.Range("D43").Value = "'4.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.12%  '

# Row 44
$ws.Range("D44").Value = "'1.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").Value = "'1.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.56%  '

# Row 46
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = "'39.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.59%  '

# Row 47
$ws.Range("D47").Value = "'0.0869"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.15%  '

# Row 48
$ws.Range("D48").Value = "'2.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.65%  '

# Row 49
$ws.Range("D49").Value = "'0.570"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.31%  '

# Row 50
$ws.Range("D50").Value = "'3.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.84%  '

# Row 51
$ws.Range("D51").Value = "'0.639"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.85%  '
